$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

$ws.Range("D2").Value = "33,33 TL - 33,33 TL"

$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""

$ws.Range("D6").Value = "8.300,01 TL - 199,41 TL"

$ws.Range("D7").Value = "%1,6"

$ws.Range("F8").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("F10").Value = ""

$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"

$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("F13").Value = ""

$ws.Range("D14").Value = "3.500 TL - 13.500 TL"
$ws.Range("F14").Value = ""
